# Fill in the next timelog entry (row 7) and update the window/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry: 2014-05-31, 18:00 -> 19:30, activity "SharpDX recherchen, SharpDX.WPF einbinden"
$ws.Range("A7").Value = 41790
$ws.Range("B7").Value = 0.75
$ws.Range("C7").Value = 0.8125
$ws.Range("E7").Value = "SharpDX recherchen, SharpDX.WPF einbinden"

# Move the active selection from A6 to A7
$ws.Range("A7").Select()

# Shift the workbook window to the right (xWindow 7020 -> 8715)
$wb.Windows.Item(1).Left = 8715
